# Fruta / hortaliza, semanal
# Insert two new weekly price rows for Membrillo (Vega Central Mapocho de
# Santiago) above the existing data, pushing the previous rows 18-50 down
# to rows 20-52.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 18:50 down by two rows.
$ws.Rows("18:19").Insert()

# New row 18
$ws.Cells.Item(18, 1).Value = 9
$ws.Cells.Item(18, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(18, 3).Value = "Metropolitana"
$ws.Cells.Item(18, 4).Value = 45002
$ws.Cells.Item(18, 5).Value = 13
$ws.Cells.Item(18, 6).Value = "Fruta"
$ws.Cells.Item(18, 7).Value = 100104
$ws.Cells.Item(18, 8).Value = "Frutos de pepita"
$ws.Cells.Item(18, 9).Value = 100104003
$ws.Cells.Item(18, 10).Value = "Membrillo"
$ws.Cells.Item(18, 11).Value = "Champion"
$ws.Cells.Item(18, 12).Value = "Especial"
$ws.Cells.Item(18, 13).Value = 45
$ws.Cells.Item(18, 14).Value = 11000
$ws.Cells.Item(18, 15).Value = 11000
$ws.Cells.Item(18, 16).Value = 11000
$ws.Cells.Item(18, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(18, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(18, 19).Value = 611
$ws.Cells.Item(18, 20).Value = 18

# New row 19
$ws.Cells.Item(19, 1).Value = 9
$ws.Cells.Item(19, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(19, 3).Value = "Metropolitana"
$ws.Cells.Item(19, 4).Value = 45002
$ws.Cells.Item(19, 5).Value = 13
$ws.Cells.Item(19, 6).Value = "Fruta"
$ws.Cells.Item(19, 7).Value = 100104
$ws.Cells.Item(19, 8).Value = "Frutos de pepita"
$ws.Cells.Item(19, 9).Value = 100104003
$ws.Cells.Item(19, 10).Value = "Membrillo"
$ws.Cells.Item(19, 11).Value = "Champion"
$ws.Cells.Item(19, 12).Value = "Primera"
$ws.Cells.Item(19, 13).Value = 80
$ws.Cells.Item(19, 14).Value = 9000
$ws.Cells.Item(19, 15).Value = 9000
$ws.Cells.Item(19, 16).Value = 9000
$ws.Cells.Item(19, 17).Value = "$/caja 18 kilos granel"
$ws.Cells.Item(19, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(19, 19).Value = 500
$ws.Cells.Item(19, 20).Value = 18
